# Apply "Multiple: Add Pre-1922 rights and more detail on how inflow split"
# edits to the Pre/Post Compact Water Rights workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MinValues")
$ws2 = $wb.Worksheets.Item("MaxValues")

# ---------------------------------------------------------------------
# MinValues sheet (sheet1)
# ---------------------------------------------------------------------

# Colorado's total (D4) is now a directly-entered value instead of the
# sum of B4+C4, and Wyoming's share (C4) is now derived from the new
# total minus the known B4 share.
$ws1.Range("D4").Value = 0.98
$ws1.Range("C4").Formula = "=D4-B4"

# New Mexico (row 7) reference source text changes from "USBR, 2020" to
# "USBR (2020)".
$ws1.Range("E7").Value = "USBR (2020)"

# New Mexico's post-1922 total (D7) is updated and now uses the same
# number format as the other post-1922 totals (D4/D5/D6 style).
$ws1.Range("D5").Copy() | Out-Null
$ws1.Range("D7").PasteSpecial(-4122) | Out-Null
$ws1.Range("D7").Value = 0.415

# ---------------------------------------------------------------------
# MaxValues sheet (sheet2)
# ---------------------------------------------------------------------

# Same restructuring of Colorado's row as on MinValues.
$ws2.Range("D4").Value = 0.983
$ws2.Range("C4").Formula = "=D4-B4"

# New Mexico's post-1922 total (D7) updated value + number format,
# and the (previously blank) source cell now cites USBR (2020).
$ws2.Range("D5").Copy() | Out-Null
$ws2.Range("D7").PasteSpecial(-4122) | Out-Null
$ws2.Range("D7").Value = 0.415
$ws2.Range("E7").Value = "USBR (2020)"

# Add the full USBR citation to the References section.
$ws2.Range("A13").Value = "USBR. (2020). `"Upper Colorado River Basin Consumptive Uses and Losses 2016 – 2020 Upper Colorado Region `", Department of Interior. https://www.usbr.gov/uc/DocLibrary/Reports/ConsumptiveUsesLosses/20210800-ProvisionalUpperColoradoRiverBasin2016-2020-CULReport-508-UCRO.pdf."

# Fill in Pre/Post-1922 values for New Mexico (MinValues sheet) with "??"
# placeholders, matching the general-format/centered/bordered style used
# elsewhere in the table (same formatting as the blank cell at C10).
$ws1.Range("C10").Copy() | Out-Null
$ws1.Range("B7").PasteSpecial(-4122) | Out-Null
$ws1.Range("C7").PasteSpecial(-4122) | Out-Null
$ws1.Range("B7").Value = "??"
$ws1.Range("C7").Value = "??"

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------

# MaxValues keeps a selection on D7 but is no longer the active tab.
$ws2.Activate() | Out-Null
$ws2.Range("D7").Select() | Out-Null

# MinValues becomes the active sheet, selected on the new B7:C7 cells.
$ws1.Activate() | Out-Null
$ws1.Range("B7:C7").Select() | Out-Null
